# Update the EC (Estado de Cuenta) database: the "Periodo Mora"/"Valor Mora"
# table (rows 16-28, columns E/F) is re-sorted into ascending period order.
# Each (period, value) pair keeps travelling together, so after the sort:
#   2110/2111/2112/2201..2209 -> 36400 ; 2209/2210 stay at 40000 (now last)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(16, 5).Value = "2110"
$ws.Cells.Item(16, 6).Value = 36400

$ws.Cells.Item(17, 5).Value = "2111"
$ws.Cells.Item(17, 6).Value = 36400

$ws.Cells.Item(18, 5).Value = "2112"
$ws.Cells.Item(18, 6).Value = 36400

$ws.Cells.Item(19, 5).Value = "2201"
$ws.Cells.Item(19, 6).Value = 36400

$ws.Cells.Item(20, 5).Value = "2202"
$ws.Cells.Item(20, 6).Value = 36400

$ws.Cells.Item(21, 5).Value = "2203"
$ws.Cells.Item(21, 6).Value = 36400

$ws.Cells.Item(22, 5).Value = "2204"
$ws.Cells.Item(22, 6).Value = 36400

$ws.Cells.Item(23, 5).Value = "2205"
$ws.Cells.Item(23, 6).Value = 36400

$ws.Cells.Item(24, 5).Value = "2206"
$ws.Cells.Item(24, 6).Value = 36400

$ws.Cells.Item(25, 5).Value = "2207"
$ws.Cells.Item(25, 6).Value = 36400

$ws.Cells.Item(26, 5).Value = "2208"
$ws.Cells.Item(26, 6).Value = 36400

$ws.Cells.Item(27, 5).Value = "2209"
$ws.Cells.Item(27, 6).Value = 40000

$ws.Cells.Item(28, 5).Value = "2210"
$ws.Cells.Item(28, 6).Value = 40000
